$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/18/2023  Through  12/24/2023"

# --- Weekly crime statistics updates ---
# Row 14
$ws.Range("J14").Value = 61
$ws.Range("K14").Value = -3.278688524590
$ws.Range("L14").Value = 9.259259259259
$ws.Range("M14").Value = -32.183908045977
$ws.Range("N14").Value = -75.619834710743
# Row 15
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 300
$ws.Range("F15").Value = 21
$ws.Range("G15").Value = 14
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 219
$ws.Range("J15").Value = 210
$ws.Range("K15").Value = 4.285714285714
$ws.Range("L15").Value = 2.816901408450
$ws.Range("M15").Value = 27.325581395348
$ws.Range("N15").Value = -59.444444444444
# Row 16
$ws.Range("C16").Value = 37
$ws.Range("D16").Value = 30
$ws.Range("E16").Value = 23.333333333333
$ws.Range("F16").Value = 134
$ws.Range("G16").Value = 144
$ws.Range("H16").Value = -6.944444444444
$ws.Range("I16").Value = 1746
$ws.Range("J16").Value = 1919
$ws.Range("K16").Value = -9.015112037519
$ws.Range("L16").Value = 22.957746478873
$ws.Range("M16").Value = -38.865546218487
$ws.Range("N16").Value = -87.326703926834
# Row 17
$ws.Range("C17").Value = 67
$ws.Range("D17").Value = 62
$ws.Range("E17").Value = 8.064516129032
$ws.Range("F17").Value = 252
$ws.Range("G17").Value = 228
$ws.Range("H17").Value = 10.526315789473
$ws.Range("I17").Value = 3462
$ws.Range("J17").Value = 3354
$ws.Range("K17").Value = 3.220035778175
$ws.Range("L17").Value = 12.548764629388
$ws.Range("M17").Value = 42.586490939044
$ws.Range("N17").Value = -47.289890377588
# Row 18
$ws.Range("C18").Value = 30
$ws.Range("D18").Value = 35
$ws.Range("E18").Value = -14.285714285714
$ws.Range("F18").Value = 114
$ws.Range("G18").Value = 165
$ws.Range("H18").Value = -30.909090909090
$ws.Range("I18").Value = 1651
$ws.Range("J18").Value = 2122
$ws.Range("K18").Value = -22.196041470311
$ws.Range("L18").Value = -6.193181818181
$ws.Range("M18").Value = -51.469723691945
$ws.Range("N18").Value = -91.053914928203
# Row 19
$ws.Range("C19").Value = 113
$ws.Range("D19").Value = 125
$ws.Range("E19").Value = -9.6
$ws.Range("F19").Value = 476
$ws.Range("G19").Value = 552
$ws.Range("H19").Value = -13.768115942029
$ws.Range("I19").Value = 6341
$ws.Range("J19").Value = 7060
$ws.Range("K19").Value = -10.184135977337
$ws.Range("L19").Value = 18.612046389824
$ws.Range("M19").Value = 15.458849235251
$ws.Range("N19").Value = -28.114726221516
# Row 20
$ws.Range("C20").Value = 33
$ws.Range("D20").Value = 27
$ws.Range("E20").Value = 22.222222222222
$ws.Range("F20").Value = 150
$ws.Range("H20").Value = 18.110236220472
$ws.Range("I20").Value = 1843
$ws.Range("J20").Value = 1764
$ws.Range("K20").Value = 4.478458049886
$ws.Range("L20").Value = 40.045592705167
$ws.Range("M20").Value = -1.601708489054
$ws.Range("N20").Value = -91.924458855490
# Row 21
$ws.Range("C21").Value = 284
$ws.Range("D21").Value = 281
$ws.Range("E21").Value = 1.067615658362
$ws.Range("F21").Value = 1153
$ws.Range("G21").Value = 1232
$ws.Range("H21").Value = -6.412337662337
$ws.Range("I21").Value = 15321
$ws.Range("J21").Value = 16490
$ws.Range("K21").Value = -7.089144936325
$ws.Range("L21").Value = 16.200227531285
$ws.Range("M21").Value = -6.063764561618
$ws.Range("N21").Value = -78.489294489294
# Row 22
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = -40
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 21
$ws.Range("H22").Value = -61.904761904761
$ws.Range("I22").Value = 189
$ws.Range("J22").Value = 196
$ws.Range("K22").Value = -3.571428571428
$ws.Range("L22").Value = 21.153846153846
$ws.Range("M22").Value = -31.272727272727
# Row 23
$ws.Range("D23").Value = 14
$ws.Range("E23").Value = -28.571428571428
$ws.Range("F23").Value = 45
$ws.Range("G23").Value = 43
$ws.Range("H23").Value = 4.651162790697
$ws.Range("I23").Value = 518
$ws.Range("J23").Value = 555
$ws.Range("K23").Value = -6.666666666666
$ws.Range("L23").Value = 7.468879668049
$ws.Range("M23").Value = 53.709198813056
# Row 24
$ws.Range("C24").Value = 310
$ws.Range("D24").Value = 273
$ws.Range("E24").Value = 13.553113553113
$ws.Range("F24").Value = 1164
$ws.Range("G24").Value = 1353
$ws.Range("H24").Value = -13.968957871396
$ws.Range("I24").Value = 15501
$ws.Range("J24").Value = 16144
$ws.Range("K24").Value = -3.982903865213
$ws.Range("L24").Value = 28.585649108253
$ws.Range("M24").Value = 27.297363882729
# Row 25
$ws.Range("C25").Value = 128
$ws.Range("D25").Value = 88
$ws.Range("E25").Value = 45.454545454545
$ws.Range("F25").Value = 481
$ws.Range("G25").Value = 379
$ws.Range("H25").Value = 26.912928759894
$ws.Range("I25").Value = 5829
$ws.Range("J25").Value = 5428
$ws.Range("K25").Value = 7.387619749447
$ws.Range("L25").Value = 18.331303288672
$ws.Range("M25").Value = -10.776060003061
# Row 26
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 500
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = 45.454545454545
$ws.Range("I26").Value = 322
$ws.Range("J26").Value = 337
$ws.Range("K26").Value = -4.451038575667
$ws.Range("L26").Value = -4.733727810650
# Row 27
$ws.Range("D27").Value = 14
$ws.Range("E27").Value = -14.285714285714
$ws.Range("F27").Value = 36
$ws.Range("G27").Value = 52
$ws.Range("H27").Value = -30.769230769230
$ws.Range("I27").Value = 622
$ws.Range("J27").Value = 688
$ws.Range("K27").Value = -9.593023255813
$ws.Range("L27").Value = -3.265940902021
# Row 28
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 5
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 20
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F28").Value = 21
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 155
$ws.Range("J28").Value = 204
$ws.Range("K28").Value = -24.019607843137
$ws.Range("L28").Value = -27.230046948356
$ws.Range("M28").Value = -43.223443223443
$ws.Range("N28").Value = -79.791395045632
# Row 29
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 4
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = 25
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F29").Value = 17
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 183.333333333333
$ws.Range("I29").Value = 132
$ws.Range("J29").Value = 161
$ws.Range("K29").Value = -18.012422360248
$ws.Range("L29").Value = -27.071823204419
$ws.Range("M29").Value = -41.850220264317
$ws.Range("N29").Value = -80.327868852459
# Row 30
$ws.Range("C30").Value = 4
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 300
$ws.Range("F30").Value = 12
$ws.Range("G30").Value = 6
$ws.Range("I30").Value = 111
$ws.Range("J30").Value = 118
$ws.Range("K30").Value = -5.932203389830
$ws.Range("L30").Value = 44.155844155844
